$wb = $excel.ActiveWorkbook

# --- Rename sheets 2 and 3 ---------------------------------------------
$wsPaginas = $wb.Worksheets.Item(2)
$wsPaginas.Name = "Páginas"

$wsActions = $wb.Worksheets.Item(3)
$wsActions.Name = "Actions"

# --- Sheet1 ("Plan1") selection: select A4:G7 ---------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A4:G7").Select()

# --- Sheet2 ("Páginas") content ------------------------------------------
$wsPaginas.Range("A1").Value = "Página"
$wsPaginas.Range("B1").Value = "Descrição"
$wsPaginas.Range("A2").Value = "Home"
$wsPaginas.Range("B2").Value = "Essa página vai mostrar um menu com opções. Usuário logado. Dashboard."

$loPaginas = $wsPaginas.ListObjects.Add(1, $wsPaginas.Range("A1:B3"), 0, 1)
$loPaginas.Name = "Tabela1"
$loPaginas.TableStyle = "TableStyleLight1"

$wsPaginas.Columns.Item(2).AutoFit()

$wsPaginas.Range("A3").Select()

# --- Sheet3 ("Actions") content -------------------------------------------
$wsActions.Range("A1").Value = "Página"
$wsActions.Range("B1").Value = "Controller"
$wsActions.Range("C1").Value = "Action"
$wsActions.Range("A2").Value = "Home"

$loActions = $wsActions.ListObjects.Add(1, $wsActions.Range("A1:C2"), 0, 1)
$loActions.Name = "Tabela2"
$loActions.TableStyle = "TableStyleLight1"

$wsActions.StandardWidth = 20

$wsActions.Range("B5").Select()

# --- Make "Actions" the active sheet/tab ----------------------------------
$wsActions.Activate()
